$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Range("G$r")
    $v = $cell.Value2
    if ($v -ne $null -and $v -ne "") {
        $parts = $v -split ", "
        if ($parts.Count -gt 1) {
            $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
            $newVal = $rotated -join ", "
            $cell.Value = $newVal
        }
    }
}
